$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.016.53'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').Value = '1.864.06'
$ws.Range('E3').Value = '  -0.91%  '
$ws.Range('D4').Value = "'1.003"
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').Value = "'312.21"
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').Value = "'1.002"
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').Value = "'0.5121"
$ws.Range('E7').Value = '  +2.13%  '
$ws.Range('D8').Value = "'0.3876"
$ws.Range('E8').Value = '  +1.44%  '
$ws.Range('D9').Value = "'0.08378"
$ws.Range('E9').Value = '  -1.65%  '
$ws.Range('D10').Value = "'1.110"
$ws.Range('E10').Value = '  -0.54%  '
$ws.Range('D11').Value = "'41.40"
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('D12').Value = "'6.157"
$ws.Range('E12').Value = '  -1.72%  '
$ws.Range('D13').Value = '1.871.43'
$ws.Range('E13').Value = '  -0.08%  '
$ws.Range('D14').Value = "'20.44"
$ws.Range('E14').Value = '  -0.73%  '
$ws.Range('D15').Value = "'7.262"
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('E16').Value = '  -0.32%  '
$ws.Range('D17').Value = "'0.00001098"
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').Value = "'90.54"
$ws.Range('E18').Value = '  -0.75%  '
$ws.Range('D19').Value = "'0.06615"
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('D20').Value = "'17.60"
$ws.Range('E20').Value = '  -2.65%  '
$ws.Range('D21').Value = "'1.002"
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('D22').Value = "'5.991"
$ws.Range('E22').Value = '  -1.65%  '
$ws.Range('D23').Value = '28.056.91'
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('D24').Value = "'11.03"
$ws.Range('E24').Value = '  -1.43%  '
$ws.Range('D25').Value = "'2.240"
$ws.Range('E25').Value = '  -1.96%  '
$ws.Range('D26').Value = '2.078.66'
$ws.Range('E26').Value = '  -0.58%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = "'2.460"
$ws.Range('E27').Value = '  -5.07%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = "'158.09"
$ws.Range('E28').Value = '  +0.65%  '
$ws.Range('D29').Value = "'20.47"
$ws.Range('E29').Value = '  -1.13%  '
$ws.Range('D30').Value = "'124.58"
$ws.Range('E30').Value = '  -1.45%  '
$ws.Range('D31').Value = "'0.1061"
$ws.Range('E31').Value = '  +0.73%  '
$ws.Range('D32').Value = "'1.027"
$ws.Range('E32').Value = '  -2.01%  '
$ws.Range('D33').Value = "'5.860"
$ws.Range('E33').Value = '  +3.84%  '
$ws.Range('D34').Value = "'3.598"
$ws.Range('E34').Value = '  -0.26%  '
$ws.Range('D35').Value = "'9.377"
$ws.Range('E35').Value = '  -2.19%  '
$ws.Range('D36').Value = "'0.02429"
$ws.Range('E36').Value = '  -1.27%  '
$ws.Range('D37').Value = "'0.06530"
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').Value = "'0.2178"
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('D39').Value = "'1.199"
$ws.Range('E39').Value = '  -2.87%  '
$ws.Range('D40').Value = "'0.6471"
$ws.Range('E40').Value = '  +1.68%  '
$ws.Range('D41').Value = "'4.969"
$ws.Range('E41').Value = '  +1.72%  '
$ws.Range('D42').Value = "'1.215"
$ws.Range('E42').Value = '  -1.88%  '
$ws.Range('E43').Value = '  -0.53%  '
$ws.Range('D44').Value = "'0.6065"
$ws.Range('E44').Value = '  +0.62%  '
$ws.Range('D45').Value = "'12.93"
$ws.Range('E45').Value = '  -1.82%  '
$ws.Range('D46').Value = "'1.285"
$ws.Range('E46').Value = '  -1.01%  '
$ws.Range('D47').Value = "'3.669"
$ws.Range('E47').Value = '  -0.45%  '
$ws.Range('D48').Value = "'1.999"
$ws.Range('E48').Value = '  +0.19%  '
$ws.Range('D49').Value = "'1.216"
$ws.Range('E49').Value = '  -0.23%  '
$ws.Range('D50').Value = "'120.75"
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('D51').Value = "'77.90"
$ws.Range('E51').Value = '  -3.51%  '
